$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Activate()

# Insert a new row for Samoa (Oceania) right after the existing Oceania block,
# before the South America block begins (i.e. at row 256).
$ws.Rows.Item(256).Insert()
$ws.Cells.Item(256, 1).Value = "Oceania"
$ws.Cells.Item(256, 2).Value = "Samoa"
$ws.Cells.Item(256, 3).Value = "Samoa"

# Reflect the author's on-screen navigation state at save time.
$ws.Range("A257").Select() | Out-Null
